$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Expand the "表3" table by one row; this also grows the table ref and
# autofilter range (A3:H7 -> A3:H8), same as Excel does when a new row is
# entered directly below an existing table.
$lo = $ws.ListObjects.Item(1)
$newRow = $lo.ListRows.Add()

# Pick up the formatting used by the rest of the table data rows before
# filling in the new values (ListRows.Add leaves the row using the plain
# default column style).
$ws.Range("A7:H7").Copy($ws.Range("A8:H8")) | Out-Null

# New npc shop entry: id 44000005 "sdlugaoyin"
$ws.Range("A8").Value = 44000005
$ws.Range("B8").Value = "sdlugaoyin"
$ws.Range("D8").Value = "spyan;spmianfen;sphujiaofen;spyancao"
$ws.Range("C8").Value = "spyumibing;spxiangjiaonai"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0
# G8 keeps the text "true" copied from G7 above (RandomPrice column)
$ws.Range("H8").Value = 3

$ws.Range("C8").Select() | Out-Null
